# Auto-generated edit script applying scheduled-runner price/profit updates
# to the Ultros_Profits workbook sheets (one per crafting class).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 87000
$ws.Range("J3").Value = 87000
$ws.Range("L3").Value = 87000
$ws.Range("N3").Value = -87228
$ws.Range("H62").Value = 5499.75
$ws.Range("I62").Value = 3999.5
$ws.Range("J62").Value = 7000
$ws.Range("K62").Value = 3999.5
$ws.Range("L62").Value = 7000
$ws.Range("M62").Value = -3375.5
$ws.Range("N62").Value = -8248
$ws.Range("H65").Value = 5499.75
$ws.Range("I65").Value = 3999.5
$ws.Range("J65").Value = 7000
$ws.Range("K65").Value = 19997.5
$ws.Range("L65").Value = 35000
$ws.Range("M65").Value = -16877.5
$ws.Range("N65").Value = -41240
$ws.Range("H80").Value = 1440.3158
$ws.Range("J80").Value = 3316.1667
$ws.Range("L80").Value = 9948.500100000001
$ws.Range("N80").Value = -11944.5001
$ws.Range("H83").Value = 1440.3158
$ws.Range("J83").Value = 3316.1667
$ws.Range("L83").Value = 29845.5003
$ws.Range("N83").Value = -39829.5003
$ws.Range("H102").Value = 87000
$ws.Range("J102").Value = 87000
$ws.Range("L102").Value = 87000
$ws.Range("N102").Value = -93490
$ws.Range("H112").Value = 1824.75
$ws.Range("I112").Value = 1168.25
$ws.Range("K112").Value = 3504.75
$ws.Range("M112").Value = -2396.75
$ws.Range("H130").Value = 58000
$ws.Range("J130").Value = 58000
$ws.Range("L130").Value = 58000
$ws.Range("N130").Value = -68040
$ws.Range("H132").Value = 11521.741
$ws.Range("I132").Value = 1936.2084
$ws.Range("J132").Value = 57532.3
$ws.Range("K132").Value = 5808.6252
$ws.Range("L132").Value = 172596.9
$ws.Range("M132").Value = -3278.6252
$ws.Range("N132").Value = -177656.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 22686
$ws.Range("I2").Value = 33356.8
$ws.Range("J2").Value = 1344.4
$ws.Range("K2").Value = 33356.8
$ws.Range("L2").Value = 1344.4
$ws.Range("M2").Value = -33243.8
$ws.Range("N2").Value = -1570.4
$ws.Range("H26").Value = 1069.1666
$ws.Range("I26").Value = 781.4
$ws.Range("J26").Value = 2508
$ws.Range("K26").Value = 781.4
$ws.Range("L26").Value = 2508
$ws.Range("M26").Value = -451.4
$ws.Range("N26").Value = -3168
$ws.Range("H32").Value = 17248666
$ws.Range("I32").Value = 18875446
$ws.Range("K32").Value = 18875446
$ws.Range("M32").Value = -18875159
$ws.Range("H81").Value = 65000
$ws.Range("J81").Value = 65000
$ws.Range("L81").Value = 65000
$ws.Range("N81").Value = -66996
$ws.Range("H84").Value = 65000
$ws.Range("J84").Value = 65000
$ws.Range("L84").Value = 195000
$ws.Range("N84").Value = -204984
$ws.Range("H116").Value = 22686
$ws.Range("I116").Value = 33356.8
$ws.Range("J116").Value = 1344.4
$ws.Range("K116").Value = 33356.8
$ws.Range("L116").Value = 1344.4
$ws.Range("M116").Value = -31062.8
$ws.Range("N116").Value = -5932.4
$ws.Range("H130").Value = 149723.75
$ws.Range("J130").Value = 149723.75
$ws.Range("L130").Value = 149723.75
$ws.Range("N130").Value = -159763.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 22686
$ws.Range("I3").Value = 33356.8
$ws.Range("J3").Value = 1344.4
$ws.Range("K3").Value = 33356.8
$ws.Range("L3").Value = 1344.4
$ws.Range("M3").Value = -33242.8
$ws.Range("N3").Value = -1572.4
$ws.Range("H94").Value = 4492.25
$ws.Range("I94").Value = 5333.5557
$ws.Range("K94").Value = 5333.5557
$ws.Range("M94").Value = -4882.5557
$ws.Range("H134").Value = 2558.7407
$ws.Range("I134").Value = 2060.0667
$ws.Range("K134").Value = 6180.2001
$ws.Range("M134").Value = -3645.2001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1896.0952
$ws.Range("I31").Value = 1791.8055
$ws.Range("J31").Value = 2521.8333
$ws.Range("K31").Value = 1791.8055
$ws.Range("L31").Value = 2521.8333
$ws.Range("M31").Value = -1496.8055
$ws.Range("N31").Value = -3111.8333
$ws.Range("H34").Value = 1896.0952
$ws.Range("I34").Value = 1791.8055
$ws.Range("J34").Value = 2521.8333
$ws.Range("K34").Value = 1791.8055
$ws.Range("L34").Value = 2521.8333
$ws.Range("M34").Value = -1589.8055
$ws.Range("N34").Value = -2925.8333
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("H62").Value = 7499.75
$ws.Range("I62").Value = 7499.75
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 7499.75
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -6875.75
$ws.Range("H65").Value = 7499.75
$ws.Range("I65").Value = 7499.75
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 37498.75
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -34378.75
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("H114").Value = 69996.336
$ws.Range("J114").Value = 69996.336
$ws.Range("L114").Value = 69996.336
$ws.Range("N114").Value = -78674.336
$ws.Range("H132").Value = 1742.375
$ws.Range("I132").Value = 1742.375
$ws.Range("K132").Value = 5227.125
$ws.Range("M132").Value = -2697.125
$ws.Range("H134").Value = 1971.2778
$ws.Range("I134").Value = 2001.6364
$ws.Range("J134").Value = 1637.3334
$ws.Range("K134").Value = 6004.9092
$ws.Range("L134").Value = 4912.0002
$ws.Range("M134").Value = -3469.9092
$ws.Range("N134").Value = -9982.0002
$ws.Range("M44").ClearContents()
$ws.Range("N62").ClearContents()
$ws.Range("N65").ClearContents()
$ws.Range("M98").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 1314.4667
$ws.Range("J47").Value = 1939.2
$ws.Range("L47").Value = 5817.6
$ws.Range("N47").Value = -6679.6
$ws.Range("H131").Value = 3502.484
$ws.Range("J131").Value = 3748.25
$ws.Range("L131").Value = 11244.75
$ws.Range("N131").Value = -21324.75
$ws.Range("H136").Value = 3398.3845
$ws.Range("J136").Value = 4750
$ws.Range("L136").Value = 14250
$ws.Range("N136").Value = -24450
$ws.Range("H138").Value = 5913.92
$ws.Range("I138").Value = 3997.0908
$ws.Range("K138").Value = 11991.2724
$ws.Range("M138").Value = -6851.2724
$ws.Range("H139").Value = 1128
$ws.Range("I139").Value = 1128
$ws.Range("K139").Value = 3384
$ws.Range("M139").Value = 1756
$ws.Range("H140").Value = 5722.5264
$ws.Range("I140").Value = 4447.5557
$ws.Range("J140").Value = 6870
$ws.Range("K140").Value = 13342.6671
$ws.Range("L140").Value = 20610
$ws.Range("M140").Value = -8162.667099999999
$ws.Range("N140").Value = -30970

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H105").Value = 68345
$ws.Range("J105").Value = 68345
$ws.Range("L105").Value = 68345
$ws.Range("N105").Value = -75333
$ws.Range("H126").Value = 6100
$ws.Range("I126").Value = 6100
$ws.Range("K126").Value = 18300
$ws.Range("M126").Value = -15830
$ws.Range("H132").Value = 1913.7222
$ws.Range("I132").Value = 1578.3846
$ws.Range("K132").Value = 4735.1538
$ws.Range("M132").Value = -2205.1538

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3463.3
$ws.Range("J68").Value = 4495
$ws.Range("L68").Value = 4495
$ws.Range("N68").Value = -5993
$ws.Range("H71").Value = 3463.3
$ws.Range("J71").Value = 4495
$ws.Range("L71").Value = 22475
$ws.Range("N71").Value = -29963
$ws.Range("H122").Value = 5995.25
$ws.Range("I122").Value = 2663.5
$ws.Range("K122").Value = 7990.5
$ws.Range("M122").Value = -5540.5
$ws.Range("H132").Value = 3331
$ws.Range("I132").Value = 2048.3428
$ws.Range("K132").Value = 6145.028399999999
$ws.Range("M132").Value = -3615.028399999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 18522602
$ws.Range("I81").Value = 4749
$ws.Range("K81").Value = 9498
$ws.Range("M81").Value = -8437
$ws.Range("H84").Value = 18522602
$ws.Range("I84").Value = 4749
$ws.Range("K84").Value = 47490
$ws.Range("M84").Value = -42186
$ws.Range("H113").Value = 841
$ws.Range("I113").Value = 323.6
$ws.Range("J113").Value = 1875.8
$ws.Range("K113").Value = 970.8000000000001
$ws.Range("L113").Value = 5627.4
$ws.Range("M113").Value = 1199.2
$ws.Range("N113").Value = -9967.4
$ws.Range("H136").Value = 1680.2333
$ws.Range("I136").Value = 762.7917
$ws.Range("J136").Value = 5350
$ws.Range("K136").Value = 2288.3751
$ws.Range("L136").Value = 16050
$ws.Range("M136").Value = 261.6248999999998
$ws.Range("N136").Value = -21150
